$wb = $excel.ActiveWorkbook

# --- ALC sheet: 39 cell updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 5039.6816
$ws.Range("I18").Value = 1332.6666
$ws.Range("K18").Value = 1332.6666
$ws.Range("M18").Value = -1048.6666
$ws.Range("H32").Value = 1298.875
$ws.Range("J32").Value = 1413
$ws.Range("L32").Value = 1413
$ws.Range("N32").Value = -2065
$ws.Range("H100").Value = 1701.5
$ws.Range("I100").Value = 1001.3333
$ws.Range("K100").Value = 1001.3333
$ws.Range("M100").Value = -460.3333
$ws.Range("H108").Value = 68000
$ws.Range("J108").Value = 68000
$ws.Range("L108").Value = 68000
$ws.Range("N108").Value = -75680
$ws.Range("H116").Value = 13780.5
$ws.Range("J116").Value = 4925.8335
$ws.Range("L116").Value = 4925.8335
$ws.Range("N116").Value = -11809.8335
$ws.Range("H136").Value = 66021.28999999999
$ws.Range("J136").Value = 66021.28999999999
$ws.Range("L136").Value = 66021.28999999999
$ws.Range("N136").Value = -76221.28999999999
$ws.Range("H137").Value = 1204.3226
$ws.Range("I137").Value = 1024.7059
$ws.Range("K137").Value = 3074.1177
$ws.Range("M137").Value = -524.1176999999998
$ws.Range("H138").Value = 1802.6957
$ws.Range("I138").Value = 1472.2142
$ws.Range("J138").Value = 2316.7778
$ws.Range("K138").Value = 4416.642599999999
$ws.Range("L138").Value = 6950.3334
$ws.Range("M138").Value = 723.3574000000008
$ws.Range("N138").Value = -17230.3334
$ws.Range("H140").Value = 62181.8
$ws.Range("J140").Value = 62181.8
$ws.Range("L140").Value = 62181.8
$ws.Range("N140").Value = -72541.8

# --- ARM sheet: 48 cell updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3584
$ws.Range("I32").Value = 3345.1345
$ws.Range("K32").Value = 3345.1345
$ws.Range("M32").Value = -3058.1345
$ws.Range("H61").Value = 35716484
$ws.Range("I61").Value = 26317454
$ws.Range("J61").Value = 55558880
$ws.Range("K61").Value = 26317454
$ws.Range("L61").Value = 55558880
$ws.Range("M61").Value = -26317242
$ws.Range("N61").Value = -55559304
$ws.Range("H74").Value = 1136.85
$ws.Range("I74").Value = 868.5806
$ws.Range("J74").Value = 2060.889
$ws.Range("K74").Value = 868.5806
$ws.Range("L74").Value = 2060.889
$ws.Range("M74").Value = 5.419399999999996
$ws.Range("N74").Value = -3808.889
$ws.Range("H77").Value = 1136.85
$ws.Range("I77").Value = 868.5806
$ws.Range("J77").Value = 2060.889
$ws.Range("K77").Value = 4342.903
$ws.Range("L77").Value = 10304.445
$ws.Range("M77").Value = 25.09699999999975
$ws.Range("N77").Value = -19040.445
$ws.Range("H97").Value = 923.25
$ws.Range("I97").Value = 861.7273
$ws.Range("K97").Value = 861.7273
$ws.Range("M97").Value = -365.7273
$ws.Range("H126").Value = 500
$ws.Range("I126").Value = 500
$ws.Range("K126").Value = 1500
$ws.Range("M126").Value = 970
$ws.Range("H132").Value = 1416.6123
$ws.Range("I132").Value = 1104.7693
$ws.Range("K132").Value = 3314.3079
$ws.Range("M132").Value = -784.3078999999998
$ws.Range("H135").Value = 32414
$ws.Range("J135").Value = 32414
$ws.Range("L135").Value = 32414
$ws.Range("N135").Value = -42554
$ws.Range("H136").Value = 35716484
$ws.Range("I136").Value = 26317454
$ws.Range("J136").Value = 55558880
$ws.Range("K136").Value = 78952362
$ws.Range("L136").Value = 166676640
$ws.Range("M136").Value = -78949812
$ws.Range("N136").Value = -166681740

# --- BSM sheet: 19 cell updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2536.4
$ws.Range("I20").Value = 2127.7778
$ws.Range("J20").Value = 3149.3333
$ws.Range("K20").Value = 2127.7778
$ws.Range("L20").Value = 3149.3333
$ws.Range("M20").Value = -1880.7778
$ws.Range("N20").Value = -3643.3333
$ws.Range("H76").Value = 61097.8
$ws.Range("J76").Value = 61097.8
$ws.Range("L76").Value = 61097.8
$ws.Range("N76").Value = -61727.8
$ws.Range("H79").Value = 61097.8
$ws.Range("J79").Value = 61097.8
$ws.Range("L79").Value = 61097.8
$ws.Range("N79").Value = -63281.8
$ws.Range("H134").Value = 4863.1035
$ws.Range("I134").Value = 5181.148
$ws.Range("K134").Value = 15543.444
$ws.Range("M134").Value = -13008.444

# --- CRP sheet: 12 cell updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1704
$ws.Range("I99").Value = 1619.4286
$ws.Range("K99").Value = 1619.4286
$ws.Range("M99").Value = -121.4286
$ws.Range("H105").Value = 1294.5883
$ws.Range("I105").Value = 1294.5883
$ws.Range("K105").Value = 1294.5883
$ws.Range("M105").Value = 452.4117000000001
$ws.Range("H126").Value = 1704
$ws.Range("I126").Value = 1619.4286
$ws.Range("K126").Value = 4858.2858
$ws.Range("M126").Value = -2388.2858

# --- CUL sheet: 4 cell updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5384888.5
$ws.Range("J131").Value = 9298.341
$ws.Range("L131").Value = 27895.023
$ws.Range("N131").Value = -37975.023

# --- GSM sheet: 25 cell updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 4142111
$ws.Range("I7").Value = 5791666.5
$ws.Range("J7").Value = 842999.7
$ws.Range("K7").Value = 5791666.5
$ws.Range("L7").Value = 842999.7
$ws.Range("M7").Value = -5791554.5
$ws.Range("N7").Value = -843223.7
$ws.Range("H8").Value = 4142111
$ws.Range("I8").Value = 5791666.5
$ws.Range("J8").Value = 842999.7
$ws.Range("K8").Value = 5791666.5
$ws.Range("L8").Value = 842999.7
$ws.Range("M8").Value = -5791527.5
$ws.Range("N8").Value = -843277.7
$ws.Range("H57").Value = 39977
$ws.Range("J57").Value = 39977
$ws.Range("L57").Value = 39977
$ws.Range("N57").Value = -41617
$ws.Range("H102").Value = 3164.9
$ws.Range("I102").Value = 3238.7778
$ws.Range("J102").Value = 2500
$ws.Range("K102").Value = 3238.7778
$ws.Range("L102").Value = 2500
$ws.Range("M102").Value = -1616.7778
$ws.Range("N102").Value = -5744

# --- LTW sheet: 4 cell updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1611.7894
$ws.Range("I132").Value = 1131.3438
$ws.Range("K132").Value = 3394.0314
$ws.Range("M132").Value = -864.0314000000003

# --- WVR sheet: 12 cell updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 7097.5654
$ws.Range("I126").Value = 10918.909
$ws.Range("K126").Value = 32756.727
$ws.Range("M126").Value = -30286.727
$ws.Range("H132").Value = 1399.3334
$ws.Range("I132").Value = 964.4815
$ws.Range("K132").Value = 2893.4445
$ws.Range("M132").Value = -363.4445000000001
$ws.Range("H136").Value = 1505.5834
$ws.Range("I136").Value = 1376.7693
$ws.Range("K136").Value = 4130.3079
$ws.Range("M136").Value = -1580.3079

Write-Output "Applied all Tonberry Profits market-data updates"